# The document opens with a stray single-letter run "w" immediately
# before the "Article X" heading text (so the paragraph currently reads
# "wArticle X"). The edit removes that stray leading run entirely,
# leaving the heading as "Article X".
#
# The stray "w" sits in its own <w:r> (distinct run properties/object
# from the following "Article X" run), so we must delete the whole
# run/character rather than just edit text in place - using Find to
# locate the exact run and Range.Delete() removes the <w:r> cleanly
# without disturbing the following run's formatting or rsid markers.

$d = $word.ActiveDocument

# Search from the very start of the document for the literal "w".
# MatchWholeWord ensures we only match the standalone leading "w" run
# (it is immediately followed by "Article", so word-boundary matching
# keeps this from over-matching into later words like "Bylaws").
$rng = $d.Range(0, 0)
$found = $rng.Find.Execute("w", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found -and $rng.Start -eq 0 -and $rng.End -eq 1) {
    $rng.Delete()
}
